# Weekly Fruta/Hortalizas price update — Macroferia Regional de Talca, Chirimoya.
# Two new daily records are inserted at the top of the data block (rows 169-170),
# pushing the existing records (previously rows 169-200) down to rows 171-202.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 169 (inserting twice at the
# same index pushes both the original row 169 and the first inserted row down).
$ws.Rows.Item(169).Insert()
$ws.Rows.Item(169).Insert()

# New row 169: Chirimoya "Especial"
$ws.Cells.Item(169,1).Value  = 5
$ws.Cells.Item(169,2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(169,3).Value  = "Maule"
$ws.Cells.Item(169,4).Value  = 45244
$ws.Cells.Item(169,5).Value  = 7
$ws.Cells.Item(169,6).Value  = "Fruta"
$ws.Cells.Item(169,7).Value  = 100107
$ws.Cells.Item(169,8).Value  = "Otros"
$ws.Cells.Item(169,9).Value  = 100107002
$ws.Cells.Item(169,10).Value = "Chirimoya"
$ws.Cells.Item(169,11).Value = "Cultivar IV Región"
$ws.Cells.Item(169,12).Value = "Especial"
$ws.Cells.Item(169,13).Value = 230
$ws.Cells.Item(169,14).Value = 20000
$ws.Cells.Item(169,15).Value = 20000
$ws.Cells.Item(169,16).Value = 20000
$ws.Cells.Item(169,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(169,18).Value = "Provincia de Limarí"
$ws.Cells.Item(169,19).Value = 2000
$ws.Cells.Item(169,20).Value = 10

# New row 170: Chirimoya "Primera"
$ws.Cells.Item(170,1).Value  = 5
$ws.Cells.Item(170,2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(170,3).Value  = "Maule"
$ws.Cells.Item(170,4).Value  = 45244
$ws.Cells.Item(170,5).Value  = 7
$ws.Cells.Item(170,6).Value  = "Fruta"
$ws.Cells.Item(170,7).Value  = 100107
$ws.Cells.Item(170,8).Value  = "Otros"
$ws.Cells.Item(170,9).Value  = 100107002
$ws.Cells.Item(170,10).Value = "Chirimoya"
$ws.Cells.Item(170,11).Value = "Cultivar IV Región"
$ws.Cells.Item(170,12).Value = "Primera"
$ws.Cells.Item(170,13).Value = 200
$ws.Cells.Item(170,14).Value = 18000
$ws.Cells.Item(170,15).Value = 18000
$ws.Cells.Item(170,16).Value = 18000
$ws.Cells.Item(170,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(170,18).Value = "Provincia de Limarí"
$ws.Cells.Item(170,19).Value = 1800
$ws.Cells.Item(170,20).Value = 10
